# Indicadores Plano de Trabalho Junho 23
# - Update the "JUNHO 23 (*)" column header to "JUNHO 23 " (drop the "(*)" footnote marker)
# - Clear the "(*) até 20 de junho" footnote text, leaving a single blank space
# - Update the June progress figures for rows 8 and 9 (column H)
# - Move the active selection down to A13 (last edited area)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = "JUNHO 23 "
$ws.Range("A12").Value = " "

$ws.Range("H8").Value = 0.90559999999999996
$ws.Range("H9").Value = 0.55000000000000004

$ws.Range("A13").Select()
